# Auto-generated script to apply scheduled market-data refresh updates
# to the Faerie_Profits leve-profit tracking workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Cells.Item(9, 8).Value = 500515.1
$ws.Cells.Item(9, 9).Value = 600374.1
$ws.Cells.Item(9, 10).Value = 1220
$ws.Cells.Item(9, 11).Value = 600374.1
$ws.Cells.Item(9, 12).Value = 1220
$ws.Cells.Item(9, 13).Value = -600205.1
$ws.Cells.Item(9, 14).Value = -1558

# Row 125
$ws.Cells.Item(125, 8).Value = 5847.5713
$ws.Cells.Item(125, 9).Value = 2249.5
$ws.Cells.Item(125, 10).Value = 6447.25
$ws.Cells.Item(125, 11).Value = 20245.5
$ws.Cells.Item(125, 12).Value = 58025.25
$ws.Cells.Item(125, 13).Value = -17785.5
$ws.Cells.Item(125, 14).Value = -62945.25

# Row 137
$ws.Cells.Item(137, 8).Value = 4195
$ws.Cells.Item(137, 9).Value = 3536.7646
$ws.Cells.Item(137, 11).Value = 10610.2938
$ws.Cells.Item(137, 13).Value = -8060.293799999999

# Row 138
$ws.Cells.Item(138, 8).Value = 142375.47
$ws.Cells.Item(138, 9).Value = 22909.2
$ws.Cells.Item(138, 10).Value = 305284.03
$ws.Cells.Item(138, 11).Value = 68727.60000000001
$ws.Cells.Item(138, 12).Value = 915852.0900000001
$ws.Cells.Item(138, 13).Value = -63587.60000000001
$ws.Cells.Item(138, 14).Value = -926132.0900000001

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 909.2
$ws.Cells.Item(2, 9).Value = 866
$ws.Cells.Item(2, 10).Value = 974
$ws.Cells.Item(2, 11).Value = 866
$ws.Cells.Item(2, 12).Value = 974
$ws.Cells.Item(2, 13).Value = -753
$ws.Cells.Item(2, 14).Value = -1200

# Row 32
$ws.Cells.Item(32, 8).Value = 2834.2654
$ws.Cells.Item(32, 9).Value = 2720.8298
$ws.Cells.Item(32, 10).Value = 5500
$ws.Cells.Item(32, 11).Value = 2720.8298
$ws.Cells.Item(32, 12).Value = 5500
$ws.Cells.Item(32, 13).Value = -2433.8298
$ws.Cells.Item(32, 14).Value = -6074

# Row 34
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 13).Value = $null

# Row 61
$ws.Cells.Item(61, 8).Value = 4541.4136
$ws.Cells.Item(61, 9).Value = 3403.3936
$ws.Cells.Item(61, 10).Value = 9499.929
$ws.Cells.Item(61, 11).Value = 3403.3936
$ws.Cells.Item(61, 12).Value = 9499.929
$ws.Cells.Item(61, 13).Value = -3191.3936
$ws.Cells.Item(61, 14).Value = -9923.929

# Row 74
$ws.Cells.Item(74, 8).Value = 1046.34
$ws.Cells.Item(74, 9).Value = 677.73334
$ws.Cells.Item(74, 10).Value = 4363.8
$ws.Cells.Item(74, 11).Value = 677.73334
$ws.Cells.Item(74, 12).Value = 4363.8
$ws.Cells.Item(74, 13).Value = 196.26666
$ws.Cells.Item(74, 14).Value = -6111.8

# Row 77
$ws.Cells.Item(77, 8).Value = 1046.34
$ws.Cells.Item(77, 9).Value = 677.73334
$ws.Cells.Item(77, 10).Value = 4363.8
$ws.Cells.Item(77, 11).Value = 3388.6667
$ws.Cells.Item(77, 12).Value = 21819
$ws.Cells.Item(77, 13).Value = 979.3333000000002
$ws.Cells.Item(77, 14).Value = -30555

# Row 102
$ws.Cells.Item(102, 8).Value = 5408024.5
$ws.Cells.Item(102, 9).Value = 2702.7
$ws.Cells.Item(102, 10).Value = 28573688
$ws.Cells.Item(102, 11).Value = 2702.7
$ws.Cells.Item(102, 12).Value = 28573688
$ws.Cells.Item(102, 13).Value = -1080.7
$ws.Cells.Item(102, 14).Value = -28576932

# Row 116
$ws.Cells.Item(116, 8).Value = 909.2
$ws.Cells.Item(116, 9).Value = 866
$ws.Cells.Item(116, 10).Value = 974
$ws.Cells.Item(116, 11).Value = 866
$ws.Cells.Item(116, 12).Value = 974
$ws.Cells.Item(116, 13).Value = 1428
$ws.Cells.Item(116, 14).Value = -5562

# Row 132
$ws.Cells.Item(132, 8).Value = 5113.159
$ws.Cells.Item(132, 9).Value = 3261.1353
$ws.Cells.Item(132, 10).Value = 14902.429
$ws.Cells.Item(132, 11).Value = 9783.4059
$ws.Cells.Item(132, 12).Value = 44707.287
$ws.Cells.Item(132, 13).Value = -7253.4059
$ws.Cells.Item(132, 14).Value = -49767.287

# Row 136
$ws.Cells.Item(136, 8).Value = 4541.4136
$ws.Cells.Item(136, 9).Value = 3403.3936
$ws.Cells.Item(136, 10).Value = 9499.929
$ws.Cells.Item(136, 11).Value = 10210.1808
$ws.Cells.Item(136, 12).Value = 28499.787
$ws.Cells.Item(136, 13).Value = -7660.1808
$ws.Cells.Item(136, 14).Value = -33599.787

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 909.2
$ws.Cells.Item(3, 9).Value = 866
$ws.Cells.Item(3, 10).Value = 974
$ws.Cells.Item(3, 11).Value = 866
$ws.Cells.Item(3, 12).Value = 974
$ws.Cells.Item(3, 13).Value = -752
$ws.Cells.Item(3, 14).Value = -1202

# Row 22
$ws.Cells.Item(22, 8).Value = 1000
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 1000
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 1000
$ws.Cells.Item(22, 13).Value = $null
$ws.Cells.Item(22, 14).Value = -1346

# Row 29
$ws.Cells.Item(29, 8).Value = 12474
$ws.Cells.Item(29, 9).Value = 11842.5
$ws.Cells.Item(29, 11).Value = 11842.5
$ws.Cells.Item(29, 13).Value = -11553.5

# Row 36
$ws.Cells.Item(36, 8).Value = 1587.6
$ws.Cells.Item(36, 9).Value = 1587.6
$ws.Cells.Item(36, 11).Value = 1587.6
$ws.Cells.Item(36, 13).Value = -1053.6

# Row 54
$ws.Cells.Item(54, 8).Value = 11966.667
$ws.Cells.Item(54, 9).Value = 2950
$ws.Cells.Item(54, 11).Value = 2950
$ws.Cells.Item(54, 13).Value = -2466

# Row 75
$ws.Cells.Item(75, 8).Value = 13450.417
$ws.Cells.Item(75, 9).Value = 6671.2856
$ws.Cells.Item(75, 10).Value = 22941.2
$ws.Cells.Item(75, 11).Value = 6671.2856
$ws.Cells.Item(75, 12).Value = 22941.2
$ws.Cells.Item(75, 13).Value = -5735.2856
$ws.Cells.Item(75, 14).Value = -24813.2

# Row 78
$ws.Cells.Item(78, 8).Value = 13450.417
$ws.Cells.Item(78, 9).Value = 6671.2856
$ws.Cells.Item(78, 10).Value = 22941.2
$ws.Cells.Item(78, 11).Value = 20013.8568
$ws.Cells.Item(78, 12).Value = 68823.60000000001
$ws.Cells.Item(78, 13).Value = -15333.8568
$ws.Cells.Item(78, 14).Value = -78183.60000000001

# Row 82
$ws.Cells.Item(82, 8).Value = 18127.924
$ws.Cells.Item(82, 9).Value = 5239.7
$ws.Cells.Item(82, 10).Value = 61088.668
$ws.Cells.Item(82, 11).Value = 5239.7
$ws.Cells.Item(82, 12).Value = 61088.668
$ws.Cells.Item(82, 13).Value = -4856.7
$ws.Cells.Item(82, 14).Value = -61854.668

# Row 85
$ws.Cells.Item(85, 8).Value = 18127.924
$ws.Cells.Item(85, 9).Value = 5239.7
$ws.Cells.Item(85, 10).Value = 61088.668
$ws.Cells.Item(85, 11).Value = 5239.7
$ws.Cells.Item(85, 12).Value = 61088.668
$ws.Cells.Item(85, 13).Value = -3913.7
$ws.Cells.Item(85, 14).Value = -63740.668

# Row 94
$ws.Cells.Item(94, 8).Value = 1356.875
$ws.Cells.Item(94, 9).Value = 1067.7142
$ws.Cells.Item(94, 10).Value = 3381
$ws.Cells.Item(94, 11).Value = 1067.7142
$ws.Cells.Item(94, 12).Value = 3381
$ws.Cells.Item(94, 13).Value = -616.7141999999999
$ws.Cells.Item(94, 14).Value = -4283

# Row 99
$ws.Cells.Item(99, 8).Value = 4816.5
$ws.Cells.Item(99, 9).Value = 3578.4285
$ws.Cells.Item(99, 10).Value = 6549.8
$ws.Cells.Item(99, 11).Value = 3578.4285
$ws.Cells.Item(99, 12).Value = 6549.8
$ws.Cells.Item(99, 13).Value = -2080.4285
$ws.Cells.Item(99, 14).Value = -9545.799999999999

# Row 134
$ws.Cells.Item(134, 8).Value = 3460.1638
$ws.Cells.Item(134, 9).Value = 1254.3489
$ws.Cells.Item(134, 10).Value = 8729.611000000001
$ws.Cells.Item(134, 11).Value = 3763.0467
$ws.Cells.Item(134, 12).Value = 26188.833
$ws.Cells.Item(134, 13).Value = -1228.0467
$ws.Cells.Item(134, 14).Value = -31258.833

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 1948.3077
$ws.Cells.Item(16, 9).Value = 987.5625
$ws.Cells.Item(16, 11).Value = 987.5625
$ws.Cells.Item(16, 13).Value = -700.5625

# Row 31
$ws.Cells.Item(31, 8).Value = 2361.2827
$ws.Cells.Item(31, 9).Value = 1710.4286
$ws.Cells.Item(31, 10).Value = 4432.1816
$ws.Cells.Item(31, 11).Value = 1710.4286
$ws.Cells.Item(31, 12).Value = 4432.1816
$ws.Cells.Item(31, 13).Value = -1415.4286
$ws.Cells.Item(31, 14).Value = -5022.1816

# Row 34
$ws.Cells.Item(34, 8).Value = 2361.2827
$ws.Cells.Item(34, 9).Value = 1710.4286
$ws.Cells.Item(34, 10).Value = 4432.1816
$ws.Cells.Item(34, 11).Value = 1710.4286
$ws.Cells.Item(34, 12).Value = 4432.1816
$ws.Cells.Item(34, 13).Value = -1508.4286
$ws.Cells.Item(34, 14).Value = -4836.1816

# Row 58
$ws.Cells.Item(58, 8).Value = 2050.8
$ws.Cells.Item(58, 9).Value = 2050.8
$ws.Cells.Item(58, 11).Value = 2050.8
$ws.Cells.Item(58, 13).Value = -1847.8

# Row 62
$ws.Cells.Item(62, 8).Value = 4138.75
$ws.Cells.Item(62, 9).Value = 4518.3335
$ws.Cells.Item(62, 10).Value = 3000
$ws.Cells.Item(62, 11).Value = 4518.3335
$ws.Cells.Item(62, 12).Value = 3000
$ws.Cells.Item(62, 13).Value = -3894.3335
$ws.Cells.Item(62, 14).Value = -4248

# Row 65
$ws.Cells.Item(65, 8).Value = 4138.75
$ws.Cells.Item(65, 9).Value = 4518.3335
$ws.Cells.Item(65, 10).Value = 3000
$ws.Cells.Item(65, 11).Value = 22591.6675
$ws.Cells.Item(65, 12).Value = 15000
$ws.Cells.Item(65, 13).Value = -19471.6675
$ws.Cells.Item(65, 14).Value = -21240

# Row 99
$ws.Cells.Item(99, 8).Value = 5203
$ws.Cells.Item(99, 9).Value = 3667
$ws.Cells.Item(99, 10).Value = 7507
$ws.Cells.Item(99, 11).Value = 3667
$ws.Cells.Item(99, 12).Value = 7507
$ws.Cells.Item(99, 13).Value = -2169
$ws.Cells.Item(99, 14).Value = -10503

# Row 113
$ws.Cells.Item(113, 8).Value = 1948.3077
$ws.Cells.Item(113, 9).Value = 987.5625
$ws.Cells.Item(113, 11).Value = 987.5625
$ws.Cells.Item(113, 13).Value = 1182.4375

# Row 122
$ws.Cells.Item(122, 8).Value = 3152.8667
$ws.Cells.Item(122, 9).Value = 3065.5186
$ws.Cells.Item(122, 10).Value = 3939
$ws.Cells.Item(122, 11).Value = 9196.5558
$ws.Cells.Item(122, 12).Value = 11817
$ws.Cells.Item(122, 13).Value = -6746.5558
$ws.Cells.Item(122, 14).Value = -16717

# Row 126
$ws.Cells.Item(126, 8).Value = 5203
$ws.Cells.Item(126, 9).Value = 3667
$ws.Cells.Item(126, 10).Value = 7507
$ws.Cells.Item(126, 11).Value = 11001
$ws.Cells.Item(126, 12).Value = 22521
$ws.Cells.Item(126, 13).Value = -8531
$ws.Cells.Item(126, 14).Value = -27461

# Row 132
$ws.Cells.Item(132, 8).Value = 2680.5356
$ws.Cells.Item(132, 9).Value = 2627.2917
$ws.Cells.Item(132, 11).Value = 7881.875100000001
$ws.Cells.Item(132, 13).Value = -5351.875100000001

# Row 134
$ws.Cells.Item(134, 8).Value = 2778.5278
$ws.Cells.Item(134, 9).Value = 2428.0688
$ws.Cells.Item(134, 10).Value = 4230.4287
$ws.Cells.Item(134, 11).Value = 7284.2064
$ws.Cells.Item(134, 12).Value = 12691.2861
$ws.Cells.Item(134, 13).Value = -4749.2064
$ws.Cells.Item(134, 14).Value = -17761.2861

# Row 136
$ws.Cells.Item(136, 8).Value = 2050.8
$ws.Cells.Item(136, 9).Value = 2050.8
$ws.Cells.Item(136, 11).Value = 6152.400000000001
$ws.Cells.Item(136, 13).Value = -3602.400000000001

$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Cells.Item(33, 8).Value = 880.6667
$ws.Cells.Item(33, 9).Value = 850
$ws.Cells.Item(33, 10).Value = 896
$ws.Cells.Item(33, 11).Value = 5100
$ws.Cells.Item(33, 12).Value = 5376
$ws.Cells.Item(33, 13).Value = -4817
$ws.Cells.Item(33, 14).Value = -5942

$ws = $wb.Worksheets.Item("GSM")
# Row 86
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 14).Value = $null

# Row 89
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 14).Value = $null

# Row 102
$ws.Cells.Item(102, 8).Value = 40141.54
$ws.Cells.Item(102, 9).Value = 1820.1666
$ws.Cells.Item(102, 11).Value = 1820.1666
$ws.Cells.Item(102, 13).Value = -198.1666

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 28324.4
$ws.Cells.Item(22, 9).Value = 10749.25
$ws.Cells.Item(22, 10).Value = 40041.168
$ws.Cells.Item(22, 11).Value = 10749.25
$ws.Cells.Item(22, 12).Value = 40041.168
$ws.Cells.Item(22, 13).Value = -10454.25
$ws.Cells.Item(22, 14).Value = -40631.168

# Row 27
$ws.Cells.Item(27, 8).Value = 28324.4
$ws.Cells.Item(27, 9).Value = 10749.25
$ws.Cells.Item(27, 10).Value = 40041.168
$ws.Cells.Item(27, 11).Value = 10749.25
$ws.Cells.Item(27, 12).Value = 40041.168
$ws.Cells.Item(27, 13).Value = -10642.25
$ws.Cells.Item(27, 14).Value = -40255.168
